$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1651168791527479"
$ws1.Range("B2").Value = "go_stims-1651168791496966.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687915108476.csv"
$ws1.Range("B4").Value = "go_stims-16511687915118477.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687915265076.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651168794857426"
$ws2.Range("B2").Value = "OB-1651168792927475.csv"
$ws2.Range("B3").Value = "TB-16511687943787577.csv"
$ws2.Range("B4").Value = "OB-1651168793189101.csv"
$ws2.Range("B5").Value = "ZB-match_1-16511687923196144.csv"
$ws2.Range("B6").Value = "TB-16511687948398352.csv"
$ws2.Range("B7").Value = "OB-16511687935073428.csv"
$ws2.Range("B8").Value = "ZB-match_0-16511687920119143.csv"
$ws2.Range("B9").Value = "ZB-match_2-16511687920622752.csv"
$ws2.Range("B10").Value = "TB-16511687936383483.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651168794859428"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651168794903962"
$ws4.Range("B2").Value = "MM_stims-1651168794872627.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687948624885.csv"
$ws4.Range("B4").Value = "MM_stims-16511687948882263.csv"
$ws4.Range("B5").Value = "ZM_stims-1651168794873633.csv"
$ws4.Range("B6").Value = "MM_stims-1651168794903962.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687948892212.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511687949824743"
$ws5.Range("B2").Value = "vSAT_stims-16511687949668422.csv"
$ws5.Range("B3").Value = "SAT_stims-1651168794910347.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687949358118.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511687949507236.csv"
